# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The diff shows every value in C2:C84 changing from 46060 to 46061,
# i.e. the underlying date serial number is incremented by one day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine last used row based on column A (Beteckning) which is populated
# for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    $v = $cell.Value2
    if ($v -ne $null) {
        $cell.Value2 = $v + 1
    }
}
